$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.255.75"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.028.48"
$ws.Range("E3").Value = "  -1.84%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.610"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.96%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "55.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.378"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.86%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0792"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("E11").Value = "  -5.09%  "
$ws.Range("D12").Value = "2.328.61"
$ws.Range("E12").Value = "  -1.83%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.54%  "
$ws.Range("E15").Value = "  -2.59%  "
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").Value = "2.019.99"
$ws.Range("E17").Value = "  -2.41%  "
$ws.Range("D18").Value = "37.219.56"
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  +1.51%  "
$ws.Range("E25").Value = "  -6.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.42"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.31%  "
$ws.Range("E28").Value = "  -1.40%  "
$ws.Range("E29").Value = "  -0.74%  "
$ws.Range("E30").Value = "  -2.42%  "
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("E32").Value = "  -1.41%  "
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.47"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.69%  "
$ws.Range("E36").Value = "  +1.15%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  -4.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.46"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.18%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0217"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.20%  "
$ws.Range("D41").Value = "1.478.99"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "96.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.19%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "16.57"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0923"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.33%  "
$ws.Range("E46").Value = "  -4.85%  "
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.77%  "
$ws.Range("D50").Value = "2.215.30"
$ws.Range("E50").Value = "  -1.86%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.60"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -10.27%  "
